# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Every player row gets the team's season record: 68 wins, 94 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, one column past the existing "Unnamed: 28" (AC) column.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, centered,
# thin-bordered) by copying the format from an existing header cell instead
# of re-building it by hand, so the new cells reuse the same style record.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every player row.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 68
    $ws.Cells.Item($r, 31).Value = 94
    $ws.Cells.Item($r, 32).Value = 0
}
